$wb = $excel.ActiveWorkbook

# 1. Consume the next available id from Sheet1 (the "names" pool).
#    Row 1 is reserved/untouched; row 2 holds the next id to use.
$ws1 = $wb.Worksheets.Item("Sheet1")
$newId = $ws1.Cells.Item(2, 1).Value()
$ws1.Rows.Item(2).Delete()

# 2. Record the newly used id on the "used" sheet.
$ws2 = $wb.Worksheets.Item("used")
$usedRange = $ws2.UsedRange
$lastRow = $usedRange.Row() + $usedRange.Rows.Count - 1
$newRow = $lastRow + 1
$ws2.Cells.Item($newRow, 1).Value = $newId
$ws2.Cells.Item($newRow, 2).Value = "ChatGPT Image 2026年1月21日 15_24_30.png"
$ws2.Cells.Item($newRow, 3).Value = "2026-01-21 15:31:20"
